$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 200
$ws.Range("A3").Value = 300

$ws.Range("B4").Select()
